$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Formula Text" column header (H1) using the new shared string
$ws.Range("H1").Value = "Formula Text"

# FORMULATEXT formulas showing the text of each FV formula in column F
$ws.Range("H2").Formula = "=FORMULATEXT(F2)"
$ws.Range("H3:H6").Formula = "=FORMULATEXT(F3)"

# Widen the new column to fit its contents
$ws.Columns("H").ColumnWidth = 31.451822916666668

# Move the active selection from F6 to F7
$ws.Range("F7").Select()
